# Slide 1, placeholder "Подзаголовок 2" holds the author/teacher credits.
# Second paragraph reads "Преподаватель: Ф. Каримова" and the teacher's
# initial is being corrected from "Ф." to "Р." (Ф. Каримова -> Р. Каримова).
#
# In the real edit, PowerPoint split the trailing run " Ф. Каримова" into
# three runs (" ", "Р. ", "Каримова") while re-typing the initial. We
# reproduce that by targeting just the "Ф. " substring (3 characters) of the
# paragraph through the Characters() sub-range and replacing it with "Р. ",
# which naturally causes the host to split the original run the same way.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$shape = $s.Shapes.Item(2)
$tf = $shape.TextFrame
$tr = $tf.TextRange

$teacherPara = $tr.Paragraphs(2)

$fullText = $teacherPara.Text
$oldInitial = "Ф. "
$newInitial = "Р. "
$startPos = $fullText.IndexOf($oldInitial) + 1

$target = $teacherPara.Characters($startPos, $oldInitial.Length)
$target.Text = $newInitial
